$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the "Update automatically" date placeholder on every slide layout
#    from 10/04/2018 to 12/04/2018 (master -> custom layouts).
# ---------------------------------------------------------------------------
$m = $p.SlideMaster
$layouts = $m.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "10/04/2018") {
                $full = $tr.Characters(1, $tr.Length)
                $full.Text = "12/04/2018"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 1 text fixes.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# "Help Tab" callout - fix the "help screeen" typo and merge the two runs
# that spelled it into a single clean sentence.
$helpShape = $s.Shapes.Item("Speech Bubble: Rectangle with Corners Rounded 7")
$helpRange = $helpShape.TextFrame.TextRange
$helpPara2 = $helpRange.Paragraphs(2, 1)
$helpSub = $helpRange.Characters($helpPara2.Start, $helpPara2.Length)
$helpSub.Text = "Click on this tab to display the help screen"

# "Status Bar" callout - reword the description text.
$statusShape = $s.Shapes.Item("Speech Bubble: Rectangle with Corners Rounded 12")
$statusRange = $statusShape.TextFrame.TextRange
$statusPara2 = $statusRange.Paragraphs(2, 1)
$statusSub = $statusRange.Characters($statusPara2.Start, $statusPara2.Length)
$statusSub.Text = "Shows when the address book was last updated"
